# Applies the "Updated cryptos list" GitHub Actions refresh to Sheet1.
# Columns: A=rank(idx) B=Coin C=Link D=Price E=Volume(1h); all data cells are text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without letting the host
# coerce numeric-looking strings (e.g. "616.70", "1.00") into numbers, and
# without leaving the cell permanently tagged with a text number-format.
function Set-TextValue($range, $text) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "91.233.06"
$ws.Range("E2").Value = "  +3.70%  "
$ws.Range("D3").Value = "3.087.99"
$ws.Range("E3").Value = "  -0.75%  "
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.17%  "
Set-TextValue "D5" "217.77"
$ws.Range("E5").Value = "  +1.54%  "
Set-TextValue "D6" "616.70"
$ws.Range("E6").Value = "  -2.86%  "
$ws.Range("E7").Value = "  -3.98%  "
Set-TextValue "D8" "0.895"
$ws.Range("E8").Value = "  +8.06%  "
Set-TextValue "D9" "1.00"
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "3.088.28"
$ws.Range("E10").Value = "  -0.63%  "
Set-TextValue "D11" "0.666"
$ws.Range("E11").Value = "  +16.88%  "
$ws.Range("E12").Value = "  +5.79%  "
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "91.127.71"
$ws.Range("E14").Value = "  +3.71%  "
$ws.Range("E15").Value = "  -0.59%  "
Set-TextValue "D16" "32.89"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").Value = "3.659.16"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "3.079.43"
$ws.Range("E18").Value = "  -0.79%  "
Set-TextValue "D19" "3.47"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("E21").Value = "  +4.17%  "
Set-TextValue "D22" "434.10"
$ws.Range("E22").Value = "  +2.46%  "
Set-TextValue "D23" "8.45"
$ws.Range("E23").Value = "  +0.29%  "
Set-TextValue "D24" "5.13"
$ws.Range("E24").Value = "  +4.90%  "
Set-TextValue "D25" "5.59"
$ws.Range("E25").Value = "  +2.10%  "
Set-TextValue "D26" "83.90"
$ws.Range("E26").Value = "  +1.48%  "
Set-TextValue "D27" "11.77"
$ws.Range("E27").Value = "  +2.39%  "
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("E29").Value = "  -0.01%  "
Set-TextValue "D30" "0.166"
$ws.Range("E30").Value = "  +6.05%  "
$ws.Range("E31").Value = "  +0.23%  "
Set-TextValue "D32" "8.66"
$ws.Range("E32").Value = "  +6.12%  "
$ws.Range("E33").Value = "  -5.13%  "
Set-TextValue "D34" "515.53"
$ws.Range("E34").Value = "  +2.97%  "
Set-TextValue "D35" "6.95"
$ws.Range("E35").Value = "  +2.41%  "
Set-TextValue "D36" "0.139"
$ws.Range("E36").Value = "  -7.59%  "
$ws.Range("E37").Value = "  -0.39%  "
$ws.Range("E38").Value = "  +0.15%  "
Set-TextValue "D39" "22.99"
$ws.Range("E39").Value = "  +3.14%  "
Set-TextValue "D40" "22.31"
$ws.Range("E40").Value = "  +0.53%  "
$ws.Range("E41").Value = "  +0.14%  "
Set-TextValue "D43" "0.140"
$ws.Range("E43").Value = "  +2.46%  "
Set-TextValue "D44" "0.367"
$ws.Range("E44").Value = "  +0.51%  "
$ws.Range("E45").Value = "  +1.44%  "
Set-TextValue "D46" "0.0721"
$ws.Range("E46").Value = "  +10.82%  "
Set-TextValue "D47" "43.88"
$ws.Range("E47").Value = "  +0.59%  "
Set-TextValue "D48" "141.84"
$ws.Range("E48").Value = "  -3.03%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue "D49" "0.000261"
$ws.Range("E49").Value = "  +10.68%  "
$ws.Range("B50").Value = "Filecoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D50" "4.21"
$ws.Range("E50").Value = "  +6.57%  "
Set-TextValue "D51" "164.66"
$ws.Range("E51").Value = "  +1.47%  "
